# Update data.xlsx from QR tool: insert two new rows of scanned data
# into the "locations" sheet, at row 2 and (new) row 4.
#
# Quirks of this COM runtime worked around below:
#  - Calling a PS function with NAMED parameters (e.g. "-RowIndex 2") against
#    a function parameter mis-marshals the value when it is later used as a
#    COM indexer argument ("Invalid row"). Only POSITIONAL args are used.
#  - Assigning a pure-digit string (e.g. a phone number) to .Value lets
#    Excel's type-inference store it as a number and drop leading zeros, so
#    those specific cells are pre-formatted as text ("@") first.
#  - Assigning an empty string "" to .Value clears/omits the cell instead of
#    leaving an (empty) shared-string cell behind. The closest reproducible
#    behavior in this runtime is Range.Copy(destination) from another empty
#    cell, which at least leaves the cell present (instead of missing).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("locations")

# NOTE: named-parameter invocation (e.g. "-RowIndex 2") against a typed
# function parameter mis-marshals the value for COM indexer calls in this
# runtime, so this helper is called with positional arguments only.
function Set-RowValues($RowIndex, $Values, $TextCols) {
    for ($i = 0; $i -lt $Values.Length; $i++) {
        $col = $i + 1
        $cell = $ws.Cells.Item($RowIndex, $col)
        if ($TextCols -contains $col) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $Values[$i]
    }
}

# --- Insert new row 2: "Ca nhan Tran Thi C" record ---
$ws.Rows.Item(2).Insert(-4121)  # xlShiftDown

$row2Values = @(
    "o1il1v7hagg",
    "kf3os4gb",
    "Cá nhân Trần Thị C",
    "Madam Thu Bakery, 21C, Võ Văn Tần, Ninh Kiều, Ninh Kiều District, Cần Thơ, 94111, Vietnam",
    "https://www.google.com/maps/search/?api=1&query=10.032100,105.786400",
    "2025-08-13T11:30:41.187Z",
    "Nhà mặt tiền",
    "0123456789",
    "CN Cần Thơ II",
    "123456789012",
    "189000000000",
    "Nguyễn Văn B",
    "199b06c9666112e3",
    "05a2ba5a8101d0cf14b11a9a0ccc45b3396eb8b6ba08821e6d600aef76960b8d"
)
# phone(8), cccd(10), customerCode(11) are all-digit -> force text so
# leading zeros / full precision survive
Set-RowValues 2 $row2Values @(8, 10, 11)

# --- Insert new row 4: "Ho kinh doanh Test 3WEL" record ---
# (original row "a4do63e3aba" now sits at row 3; new row goes below it)
$ws.Rows.Item(4).Insert(-4121)  # xlShiftDown

$row4Values = @(
    "xpodt5txnr",
    "TESTM8BT",
    "Hộ kinh doanh Test 3WEL",
    "02 Hòa Bình, Ninh Kiều, Cần Thơ",
    "https://www.google.com/maps/search/?api=1&query=02%20H%C3%B2a%20B%C3%ACnh%2C%20Ninh%20Ki%E1%BB%81u%2C%20C%E1%BA%A7n%20Th%C6%A1",
    "2025-08-13T08:09:42.524Z",
    "Dòng test thêm nhanh",
    "0905167266",
    "CN Cần Thơ II",
    "0342835098143",
    "KH5912",
    "Demo User"
)
Set-RowValues 4 $row4Values @(8, 10)

# pinSalt/pinHash (cols 13/14) are blank for this record. A direct
# .Value = "" clears/omits the cell entirely in this runtime, so instead
# copy an already-blank pinSalt/pinHash cell (row 3, shifted original
# "a4do63e3aba" record) onto row 4's cells to leave them present-but-empty.
$ws.Cells.Item(3, 13).Copy($ws.Cells.Item(4, 13))
$ws.Cells.Item(3, 14).Copy($ws.Cells.Item(4, 14))

$wb.Save()
